$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K (strikeout) column values replacing old Strike# estimates
$kValues = @{
    2 = 2
    3 = 2
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 3
    18 = 0
    19 = 3
    20 = 0
    21 = 2
    22 = 0
    23 = 2
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 3
    36 = 1
    37 = 0
    38 = 1
    39 = 1
    40 = 2
    41 = 0
    42 = 2
    43 = 1
    44 = 0
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 2
    50 = 1
    51 = 0
    52 = 0
    53 = 2
    54 = 3
    55 = 0
    56 = 2
    57 = 0
    58 = 2
    59 = 0
    60 = 0
    61 = 0
    62 = 1
    63 = 1
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 3
    69 = 2
    70 = 1
    71 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}
